$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data to sheet1 (rows 2-51).
# Column D (Price) values are plain text (often using "." as thousands
# separator), so force text format before/after assignment to avoid Excel
# auto-converting them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.015.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.328.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.47%  "
$ws.Range("E7").Value = "  -3.84%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.321.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.621"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.856.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("E17").Value = "  -3.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.318.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.890.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.971"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "427.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.36%  "
$ws.Range("E24").Value = "  -2.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.48%  "
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "594.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -8.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0747"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.364"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.090.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  -5.00%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0405"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.81%  "
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.46%  "
